# Remove the "SuccessVar" entry (row 10) from the codebook worksheet.
# Deleting the whole row shifts rows 11-12 up to become rows 10-11,
# which matches the rest of the table (their "Number" column values
# were already 10 and 11, so no renumbering is needed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(10).Delete()

$ws.Rows.Item(10).Select()
